$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsDCpUC = $wb.Worksheets.Item("DCpUC")

# Add the new India:US GDP per capita adjustment factor to the About sheet
$wsAbout.Range("A32").Value = "India:US GDP per capita adjustment, see InputData/scaling-factors.xlsx"
$wsAbout.Range("A33").Value = 0.032347480211350491

# Scale every cost formula in the DCpUC sheet by the new adjustment factor
$wsDCpUC.Range("B2").Formula = '=Data!C3*About!$A$33'
$wsDCpUC.Range("B3").Formula = '=Data!C8*About!$A$33'
$wsDCpUC.Range("B4").Formula = '=Data!C13*About!$A$33'
$wsDCpUC.Range("B5").Formula = '=Data!B19*About!$A$33'
$wsDCpUC.Range("B6").Formula = '=Data!C6*About!$A$33'
$wsDCpUC.Range("B7").Formula = '=Data!C5*About!$A$33'
$wsDCpUC.Range("B8").Formula = '=Data!C4*About!$A$33'
$wsDCpUC.Range("B9").Formula = '=Data!C3*About!$A$33'
$wsDCpUC.Range("B10").Formula = '=B5*About!$A$33'
$wsDCpUC.Range("B11").Formula = '=Data!C7*About!$A$33'
$wsDCpUC.Range("B12").Formula = '=Data!C8*About!$A$33'
$wsDCpUC.Range("B13").Formula = '=Data!C3*About!$A$33'
$wsDCpUC.Range("B14").Formula = '=Data!C2*About!$A$33'
$wsDCpUC.Range("B15").Formula = '=B11*About!$A$33'
$wsDCpUC.Range("B16").Formula = '=B11*About!$A$33'
$wsDCpUC.Range("B17").Formula = '=B9*About!$A$33'

# Restore the selection/active-cell state recorded in the workbook for each sheet
$wsDCpUC.Activate()
$wsDCpUC.Range("B18").Select()

$wsAbout.Activate()
$wsAbout.Range("A32:A33").Select()
